# InfectionModels.xlsx — split the single "VC" row into two rows:
# "VC (no COE)" and "VC (COE)", each with a full set of model-fit values,
# and widen column K to fit the longer "beta pm se" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert a new blank row at 7 so the existing blank separator
# row (previously row 7) and all the rows below it (survival, MDR,
# body size, growth) shift down by one.
$ws.Rows("7:7").Insert()

# Row 6: "VC (no COE)" — fill in the previously-empty B:M cells.
$ws.Range("A6").Value = "VC (no COE)"
$ws.Range("B6").Value = "2"
$ws.Range("C6").Value = "0.445"
$ws.Range("D6").Value = "0.801"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "35.844"
$ws.Range("G6").Value = "<0.001"
$ws.Range("H6").Value = "2"
$ws.Range("I6").Value = "3.548"
$ws.Range("J6").Value = "0.1696"
$ws.Range("K6").Value = "3.803 (0.449)"
$ws.Range("L6").Value = "8.474"
$ws.Range("M6").Value = "<0.001"

# Row 7 (new): "VC (COE)"
$ws.Range("A7").Value = "VC (COE)"
$ws.Range("B7").Value = "2"
$ws.Range("C7").Value = "0.03"
$ws.Range("D7").Value = "0.979"
$ws.Range("E7").Value = "1"
$ws.Range("F7").Value = "1.172"
$ws.Range("G7").Value = "0.280"
$ws.Range("H7").Value = "2"
$ws.Range("I7").Value = "0.288"
$ws.Range("J7").Value = "0.866"
$ws.Range("K7").Value = "0.779 (0.168)"
$ws.Range("L7").Value = "4.648"
$ws.Range("M7").Value = "<0.001"

# Column K needs to be wide enough for the longer "x.xxx (x.xxx)" strings.
$ws.Range("K1").EntireColumn.ColumnWidth = 27.1640625

# Restore the cursor to where the author left it.
$ws.Range("F15").Select()
